$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: rename "The Infestor" -> "Infestor"
$ws.Range("A11").Value = "Infestor"

# Row 12: bounty (H) 6 -> 4
$ws.Range("H12").Value = 4

# Row 14: dmg-min (W) 30 -> 37, dmg-max (X) 33 -> 40
$ws.Range("W14").Value = 37
$ws.Range("X14").Value = 40

# Row 16: hp (P) 605 -> 695
$ws.Range("P16").Value = 695

# Row 21: rename "Fatty" -> "Blightbringer"
$ws.Range("A21").Value = "Blightbringer"

# Row 31: rename "Experimental Hybrid" -> "Phase Shifter"
$ws.Range("A31").Value = "Phase Shifter"
